# Daily attendance processing - 2026-01-04 07:11:33
# Reorders the "Recorded By" (column G) values so that the first two
# comma-separated entries are swapped (the second entry moves to the front).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value2 = [string]::Join(", ", $parts)
        }
    }
}
